# Fruta / hortaliza, semanal
# This edit reshuffles the per-row data for columns D, J, K, L, M, N, P, Q
# (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado,
# Unidad de comercializacion, Precio $/Kg, Kg o Unidades) across data rows
# 2-30, while columns A, B, C, E, F, G, H, I, O, R are left untouched
# (they already hold the same constant values on every row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maps destination row -> source row (source row values, as they exist in
# the *original* workbook, are written into the destination row).
$rowMap = @{
    2  = 24
    3  = 29
    4  = 23
    5  = 9
    6  = 18
    7  = 28
    8  = 3
    9  = 22
    10 = 14
    11 = 21
    12 = 15
    13 = 17
    14 = 25
    15 = 7
    16 = 8
    17 = 4
    18 = 13
    19 = 5
    20 = 10
    21 = 27
    22 = 6
    23 = 20
    24 = 26
    25 = 16
    26 = 19
    27 = 2
    28 = 30
    29 = 12
    30 = 11
}

# Columns whose values move together with a row's data.
$cols = @(4, 10, 11, 12, 13, 14, 16, 17)  # D, J, K, L, M, N, P, Q

# Snapshot the original values for every row/column first, since the
# mapping is a permutation and destination rows may also be source rows
# for other destinations.
$snapshot = @{}
foreach ($r in 2..30) {
    foreach ($c in $cols) {
        $snapshot["$r`_$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $snapshot["$srcRow`_$c"]
    }
}
